# Updates column A (match id) and column E (Round) for rows 2-27
# Column E values change from text "Matchweek N" to the plain numeric N
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  A = 1;  E = 1  },
    @{ Row = 3;  A = 4;  E = 3  },
    @{ Row = 4;  A = 7;  E = 5  },
    @{ Row = 5;  A = 9;  E = 7  },
    @{ Row = 6;  A = 11; E = 8  },
    @{ Row = 7;  A = 14; E = 10 },
    @{ Row = 8;  A = 18; E = 13 },
    @{ Row = 9;  A = 21; E = 16 },
    @{ Row = 10; A = 26; E = 19 },
    @{ Row = 11; A = 28; E = 21 },
    @{ Row = 12; A = 30; E = 23 },
    @{ Row = 13; A = 32; E = 25 },
    @{ Row = 14; A = 33; E = 22 },
    @{ Row = 15; A = 15; E = 11 },
    @{ Row = 16; A = 19; E = 14 },
    @{ Row = 17; A = 20; E = 15 },
    @{ Row = 18; A = 6;  E = 4  },
    @{ Row = 19; A = 26; E = 17 },
    @{ Row = 20; A = 17; E = 18 },
    @{ Row = 21; A = 3;  E = 2  },
    @{ Row = 22; A = 39; E = 26 },
    @{ Row = 23; A = 16; E = 12 },
    @{ Row = 24; A = 5;  E = 6  },
    @{ Row = 25; A = 8;  E = 9  },
    @{ Row = 26; A = 31; E = 24 },
    @{ Row = 27; A = 19; E = 20 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 1).Value = $u.A
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
